$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns("B").ColumnWidth = 112.16666666666667
$ws.Columns("C").ColumnWidth = 91.6

# --- Row 159: blank thick-bottom divider row ---
$ws.Rows(159).RowHeight = 15

# --- Copy formatting templates for new rows ---
$ws.Range("B4:C4").Copy()   # header style (s=1 / s=2)
$ws.Range("B160:C160").PasteSpecial(-4122)

$ws.Range("B5:C5").Copy()   # normal data-row style (s=3 / s=4)
$ws.Range("B161:C161").PasteSpecial(-4122)
$ws.Range("B162:C162").PasteSpecial(-4122)
$ws.Range("B163:C163").PasteSpecial(-4122)
$ws.Range("B164:C164").PasteSpecial(-4122)
$ws.Range("B165:C165").PasteSpecial(-4122)
$ws.Range("B166:C166").PasteSpecial(-4122)
$ws.Range("B167:C167").PasteSpecial(-4122)
$ws.Range("B168:C168").PasteSpecial(-4122)
$ws.Range("B169:C169").PasteSpecial(-4122)
$ws.Range("B170:C170").PasteSpecial(-4122)
$ws.Range("B171:C171").PasteSpecial(-4122)
$ws.Range("B172:C172").PasteSpecial(-4122)
$ws.Range("B173:C173").PasteSpecial(-4122)
$ws.Range("B174:C174").PasteSpecial(-4122)
$ws.Range("B175:C175").PasteSpecial(-4122)
$ws.Range("B176:C176").PasteSpecial(-4122)

$ws.Range("B19:C19").Copy()   # final thick-bottom data-row style (s=5 / s=6)
$ws.Range("B177:C177").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Values ---
$ws.Range("B160").Value = 'Команда'
$ws.Range("C160").Value = 'Что делает'
$ws.Range("B161").Value = 'helm version'
$ws.Range("C161").Value = 'Пока версию Helm'
$ws.Range("B162").Value = 'helm list'
$ws.Range("C162").Value = 'Показать все задеплоенные Helm Releases'
$ws.Range("B164").Value = 'helm search hub'
$ws.Range("C164").Value = 'Показать Helm Chart с общего списка Hub'
$ws.Range("B165").Value = 'helm search repo'
$ws.Range("C165").Value = 'Показать Helm Chart из добавленных Repos'
$ws.Range("B167").Value = 'helm install app1 Denis-Chart/'
$ws.Range("C167").Value = 'Задеплоить Helm Chart app1 из директории Denis-Chart'
$ws.Range("B168").Value = 'helm upgrade app1 Denis-Chart/ --set container.image=adv4000/k8sphp:version2'
$ws.Range("C168").Value = 'Обновить Деплоймент app1'
$ws.Range("B170").Value = 'helm create MyChart'
$ws.Range("C170").Value = 'Сделать скелет Helm Chart в директории MyChart'
$ws.Range("B171").Value = 'helm package Denis-Chart/'
$ws.Range("C171").Value = 'Запаковать Helm Chart в tgz архив'
$ws.Range("B172").Value = 'helm install app2 App-HelmChart-0.1.0.tgz'
$ws.Range("C172").Value = 'Задеплоить Helm Chart app2 из архива'
$ws.Range("B173").Value = 'helm delete app1'
$ws.Range("C173").Value = 'Удалить Деплоймент Helm Chart app1'
$ws.Range("B174").Value = 'helm uninstall app1'
$ws.Range("C174").Value = 'Удалить Деплоймент Helm Chart app1'
$ws.Range("B176").Value = 'helm repo add bitnami https://charts.bitnami.com/bitnami'
$ws.Range("C176").Value = 'Добавить Helm Chart Repo от bitnami'
$ws.Range("B177").Value = 'helm install my_website bitnami/apache -f my_values.yaml'
$ws.Range("C177").Value = 'Задеплоить Helm Chart bitnami/apache с нашими переменными'

# --- Row heights for data rows ---
$ws.Rows(160).RowHeight = 23.4
$ws.Rows(161).RowHeight = 23.4
$ws.Rows(162).RowHeight = 23.4
$ws.Rows(163).RowHeight = 23.4
$ws.Rows(164).RowHeight = 23.4
$ws.Rows(165).RowHeight = 23.4
$ws.Rows(166).RowHeight = 23.4
$ws.Rows(167).RowHeight = 23.4
$ws.Rows(168).RowHeight = 23.4
$ws.Rows(169).RowHeight = 23.4
$ws.Rows(170).RowHeight = 23.4
$ws.Rows(171).RowHeight = 23.4
$ws.Rows(172).RowHeight = 23.4
$ws.Rows(173).RowHeight = 23.4
$ws.Rows(174).RowHeight = 23.4
$ws.Rows(175).RowHeight = 23.4
$ws.Rows(176).RowHeight = 23.4
$ws.Rows(177).RowHeight = 24

# --- Selection / view ---
$ws.Range("A175").Select()
